# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column D to Text format first so numeric-looking price strings
# (e.g. "239.45") are stored as text, matching the inline-string cells
# already used throughout this sheet, rather than being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '96.946.16'
$ws.Range('E2').Value = '  -0.15%  '

$ws.Range('D3').Value = '3.670.39'
$ws.Range('E3').Value = '  +2.78%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '239.45'
$ws.Range('E5').Value = '  -0.69%  '

$ws.Range('E6').Value = '  +10.48%  '

$ws.Range('D7').Value = '654.04'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('E8').Value = '  -0.28%  '

$ws.Range('E9').Value = '  +3.47%  '

$ws.Range('E10').Value = '  +0.04%  '

$ws.Range('D11').Value = '3.668.84'
$ws.Range('E11').Value = '  +2.82%  '

$ws.Range('D12').Value = '45.50'
$ws.Range('E12').Value = '  +3.07%  '

$ws.Range('E13').Value = '  +1.24%  '

$ws.Range('D14').Value = '6.88'
$ws.Range('E14').Value = '  +7.17%  '

$ws.Range('D15').Value = '4.354.28'
$ws.Range('E15').Value = '  +2.87%  '

$ws.Range('E16').Value = '  +3.61%  '

$ws.Range('D17').Value = '96.638.31'
$ws.Range('E17').Value = '  -0.19%  '

$ws.Range('D18').Value = '9.03'
$ws.Range('E18').Value = '  +4.66%  '

$ws.Range('D19').Value = '3.686.84'
$ws.Range('E19').Value = '  +3.52%  '

$ws.Range('D20').Value = '19.01'
$ws.Range('E20').Value = '  +5.80%  '

$ws.Range('D21').Value = '12.78'
$ws.Range('E21').Value = '  +0.57%  '

$ws.Range('D22').Value = '0.533'
$ws.Range('E22').Value = '  +0.75%  '

$ws.Range('D23').Value = '531.58'
$ws.Range('E23').Value = '  +3.56%  '

$ws.Range('E24').Value = '  +0.67%  '

$ws.Range('D25').Value = '7.17'
$ws.Range('E25').Value = '  +4.67%  '

$ws.Range('E26').Value = '  -0.69%  '

$ws.Range('D27').Value = '102.55'
$ws.Range('E27').Value = '  +0.97%  '

$ws.Range('D28').Value = '13.49'
$ws.Range('E28').Value = '  +3.51%  '

$ws.Range('D29').Value = '3.867.23'
$ws.Range('E29').Value = '  +2.80%  '

$ws.Range('E30').Value = '  +0.27%  '

$ws.Range('D31').Value = '12.53'
$ws.Range('E31').Value = '  +5.52%  '

$ws.Range('D32').Value = '3.05'
$ws.Range('E32').Value = '  +1.94%  '

$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  +15.31%  '

$ws.Range('E35').Value = '  +1.44%  '

$ws.Range('D36').Value = '32.71'
$ws.Range('E36').Value = '  +3.05%  '

$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  +0.00%  '

$ws.Range('D38').Value = '656.22'
$ws.Range('E38').Value = '  +6.53%  '

$ws.Range('D39').Value = '0.607'
$ws.Range('E39').Value = '  +7.80%  '

$ws.Range('D40').Value = '8.98'
$ws.Range('E40').Value = '  +1.33%  '

$ws.Range('D41').Value = '6.96'
$ws.Range('E41').Value = '  +15.77%  '

$ws.Range('D42').Value = '0.162'
$ws.Range('E42').Value = '  +5.23%  '

$ws.Range('E43').Value = '  +2.39%  '

$ws.Range('D44').Value = '0.963'
$ws.Range('E44').Value = '  +4.51%  '

$ws.Range('D45').Value = '38.46'
$ws.Range('E45').Value = '  +16.55%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.452'
$ws.Range('E47').Value = '  +11.05%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0462'
$ws.Range('E48').Value = '  +5.48%  '

$ws.Range('E49').Value = '  +1.38%  '

$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.77'
$ws.Range('E50').Value = '  +3.04%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '23.64'
$ws.Range('E51').Value = '  +0.20%  '

# Restore the default (unstyled) cell style on column D now that the
# values are safely stored as text, so no residual number-format style
# is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
